$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.065879106521606
$ws.Range("B1").Value = 2.053476810455322
$ws.Range("C1").Value = 8.531377792358398
$ws.Range("D1").Value = 1.37977409362793
$ws.Range("E1").Value = 0.879473865032196
